# Auto-generated edit script applying scheduled market-data refresh
# to the Tonberry_Profits workbook (per commit diff).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 33375.145
$ws.Range("I21").Value = 33375.145
$ws.Range("K21").Value = 33375.145
$ws.Range("M21").Value = -32907.145

$ws.Range("H23").Value = 33375.145
$ws.Range("I23").Value = 33375.145
$ws.Range("K23").Value = 33375.145
$ws.Range("M23").Value = -33141.145

$ws.Range("H29").Value = 3633.3333
$ws.Range("I29").Value = 500
$ws.Range("K29").Value = 1500
$ws.Range("M29").Value = -1219

$ws.Range("H38").Value = 1346
$ws.Range("I38").Value = 1346
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 4038
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H51").Value = 5916.6665
$ws.Range("J51").Value = 6100
$ws.Range("L51").Value = 6100
$ws.Range("N51").Value = -7068

$ws.Range("H53").Value = 295.8889
$ws.Range("I53").Value = 259.25
$ws.Range("J53").Value = 325.2
$ws.Range("K53").Value = 259.25
$ws.Range("L53").Value = 325.2
$ws.Range("M53").Value = 377.75
$ws.Range("N53").Value = -1599.2

$ws.Range("H55").Value = 304.85715
$ws.Range("I55").Value = 244
$ws.Range("J55").Value = 365.7143
$ws.Range("K55").Value = 244
$ws.Range("L55").Value = 365.7143
$ws.Range("M55").Value = -30
$ws.Range("N55").Value = -793.7143

$ws.Range("H58").Value = 1493.1333
$ws.Range("I58").Value = 409.83334
$ws.Range("J58").Value = 2215.3333
$ws.Range("K58").Value = 1229.50002
$ws.Range("L58").Value = 6645.999899999999
$ws.Range("M58").Value = -1079.50002
$ws.Range("N58").Value = -6945.999899999999

$ws.Range("H87").Value = 57925
$ws.Range("J87").Value = 57925
$ws.Range("L87").Value = 57925
$ws.Range("N87").Value = -60421

$ws.Range("H90").Value = 57925
$ws.Range("J90").Value = 57925
$ws.Range("L90").Value = 173775
$ws.Range("N90").Value = -186255

$ws.Range("H138").Value = 1833.1522
$ws.Range("I138").Value = 1063.138
$ws.Range("K138").Value = 3189.414
$ws.Range("M138").Value = 1950.586

$ws.Range("H141").Value = 2804627.5
$ws.Range("I141").Value = 7003048
$ws.Range("K141").Value = 21009144
$ws.Range("M141").Value = -21003964

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3514.2114
$ws.Range("I32").Value = 2219.9556
$ws.Range("J32").Value = 11834.429
$ws.Range("K32").Value = 2219.9556
$ws.Range("L32").Value = 11834.429
$ws.Range("M32").Value = -1932.9556
$ws.Range("N32").Value = -12408.429

$ws.Range("H88").Value = 23769.6
$ws.Range("I88").Value = 2949.5
$ws.Range("J88").Value = 28974.625
$ws.Range("K88").Value = 2949.5
$ws.Range("L88").Value = 28974.625
$ws.Range("M88").Value = -2543.5
$ws.Range("N88").Value = -29786.625

$ws.Range("H91").Value = 23769.6
$ws.Range("I91").Value = 2949.5
$ws.Range("J91").Value = 28974.625
$ws.Range("K91").Value = 2949.5
$ws.Range("L91").Value = 28974.625
$ws.Range("M91").Value = -1545.5
$ws.Range("N91").Value = -31782.625

$ws.Range("H97").Value = 661.5
$ws.Range("I97").Value = 661.5
$ws.Range("K97").Value = 661.5
$ws.Range("M97").Value = -165.5

$ws.Range("H109").Value = 51187.75
$ws.Range("J109").Value = 51187.75
$ws.Range("L109").Value = 51187.75
$ws.Range("N109").Value = -53961.75

$ws.Range("H122").Value = 2527.7856
$ws.Range("I122").Value = 2449.0833
$ws.Range("K122").Value = 7347.249899999999
$ws.Range("M122").Value = -4897.249899999999

$ws.Range("H132").Value = 1401.119
$ws.Range("I132").Value = 947.1739
$ws.Range("J132").Value = 1950.6316
$ws.Range("K132").Value = 2841.5217
$ws.Range("L132").Value = 5851.8948
$ws.Range("M132").Value = -311.5217000000002
$ws.Range("N132").Value = -10911.8948

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 41032.75
$ws.Range("J141").Value = 41032.75
$ws.Range("L141").Value = 41032.75
$ws.Range("N141").Value = -51392.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 182.5
$ws.Range("J2").Value = 100
$ws.Range("L2").Value = 600
$ws.Range("N2").Value = -826

$ws.Range("H12").Value = 89.59999999999999
$ws.Range("I12").Value = 50
$ws.Range("J12").Value = 99.5
$ws.Range("K12").Value = 150
$ws.Range("L12").Value = 298.5
$ws.Range("M12").Value = 23
$ws.Range("N12").Value = -644.5

$ws.Range("H33").Value = 42
$ws.Range("I33").Value = 55
$ws.Range("J33").Value = 29
$ws.Range("K33").Value = 330
$ws.Range("L33").Value = 174
$ws.Range("M33").Value = -47
$ws.Range("N33").Value = -740

$ws.Range("H129").Value = 41302.277
$ws.Range("J129").Value = 56827.383
$ws.Range("L129").Value = 170482.149
$ws.Range("N129").Value = -180482.149

$ws.Range("H131").Value = 10918.072
$ws.Range("J131").Value = 12076.161
$ws.Range("L131").Value = 36228.483
$ws.Range("N131").Value = -46308.483

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3389.4167
$ws.Range("I102").Value = 4109.8335
$ws.Range("J102").Value = 2669
$ws.Range("K102").Value = 4109.8335
$ws.Range("L102").Value = 2669
$ws.Range("M102").Value = -2487.8335
$ws.Range("N102").Value = -5913

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6487.769
$ws.Range("I7").Value = 3851
$ws.Range("J7").Value = 7659.6665
$ws.Range("K7").Value = 3851
$ws.Range("L7").Value = 7659.6665
$ws.Range("M7").Value = -3739
$ws.Range("N7").Value = -7883.6665

$ws.Range("H40").Value = 6712.25
$ws.Range("I40").Value = 1833.3334
$ws.Range("J40").Value = 9639.6
$ws.Range("K40").Value = 1833.3334
$ws.Range("L40").Value = 9639.6
$ws.Range("M40").Value = -1697.3334
$ws.Range("N40").Value = -9911.6

$ws.Range("H126").Value = 6487.769
$ws.Range("I126").Value = 3851
$ws.Range("J126").Value = 7659.6665
$ws.Range("K126").Value = 11553
$ws.Range("L126").Value = 22978.9995
$ws.Range("M126").Value = -9083
$ws.Range("N126").Value = -27918.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 446.7857
$ws.Range("I100").Value = 276.7
$ws.Range("J100").Value = 872
$ws.Range("K100").Value = 553.4
$ws.Range("L100").Value = 1744
$ws.Range("M100").Value = -12.39999999999998
$ws.Range("N100").Value = -2826

$ws.Range("H107").Value = 625.95654
$ws.Range("I107").Value = 475.64285
$ws.Range("J107").Value = 859.7778
$ws.Range("K107").Value = 1426.92855
$ws.Range("L107").Value = 2579.3334
$ws.Range("M107").Value = 493.0714499999999
$ws.Range("N107").Value = -6419.3334

$ws.Range("H123").Value = 45663.156
$ws.Range("I123").Value = 30000
$ws.Range("K123").Value = 30000
$ws.Range("M123").Value = -25100
